$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (locked cells, password-hashed) so the
# protected-range figures can't be written directly. Unprotect, make the
# edits, then re-apply the same protection flags before saving so the
# workbook round-trips as "still protected" like the source file.
$ws.Unprotect("lido")

# Disclaimer banner in A11: bump the "as of" date by one day.
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-20 for illustrative purposes only and are subject to change."
$ws.Range("A11").Value = $newText

# Updated Weight (D) / Percent Change (E) figures for rows 2-8.
$ws.Range("D2").Value = 0.5013958350152756
$ws.Range("E2").Value = 0.005220338983050965

$ws.Range("D3").Value = 0.2433902084114753
$ws.Range("E3").Value = 0.01581381470983145

$ws.Range("D4").Value = 0.09498283549774354
$ws.Range("E4").Value = 0.0097361004355625

$ws.Range("D5").Value = 0.1033404628370274
$ws.Range("E5").Value = -0.0002803476310625497

$ws.Range("D6").Value = 0.03008085297552161
$ws.Range("E6").Value = -0.0001926225561013073

$ws.Range("D7").Value = 0.02680980526295662
$ws.Range("E7").Value = 0.007449084713527121

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0.007556089313354608

# Restore sheet protection with the same settings as before (contents +
# objects + scenarios locked, matching the source sheetProtection element).
$ws.Protect("lido", $false, $true, $true)
